# Add the three new names below the existing list in column A
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Pasant"
$ws.Range("A5").Value = "Suzan"
$ws.Range("A6").Value = "Johnathan"

# Auto-fit column A to the new, longer contents
$ws.Columns("A:A").AutoFit() | Out-Null

# Leave the last entered cell selected
$ws.Range("A6").Select() | Out-Null
